$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new quarter columns at D:E; existing D:K data shifts right to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting from the (now-shifted) old column D into the two new columns, per contiguous data block
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new quarter data: column D = period ending 2018-09-30, column E = period ending 2018-06-30
$ws.Range("D7").Value = 43373
$ws.Range("E7").Value = 43281
$ws.Range("D8").Value = 101200
$ws.Range("E8").Value = 111200
$ws.Range("D9").Value = 93400
$ws.Range("E9").Value = 99900
$ws.Range("D10").Value = 7800
$ws.Range("E10").Value = 11300
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 99600
$ws.Range("E17").Value = 106200
$ws.Range("D18").Value = 1600
$ws.Range("E18").Value = 5000
$ws.Range("D20").Value = 800
$ws.Range("E20").Value = -100
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 2400
$ws.Range("E23").Value = 4900
$ws.Range("D24").Value = 800
$ws.Range("E24").Value = 1500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 1700
$ws.Range("E26").Value = 3400
$ws.Range("D27").Value = 700
$ws.Range("E27").Value = 1800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -800
$ws.Range("E32").Value = 100
$ws.Range("D33").Value = 700
$ws.Range("E33").Value = 1800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 700
$ws.Range("E35").Value = 1800
$ws.Range("D38").Value = 43373
$ws.Range("E38").Value = 43281
$ws.Range("D41").Value = 48600
$ws.Range("E41").Value = 52400
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = 300
$ws.Range("D43").Value = 109800
$ws.Range("E43").Value = 117700
$ws.Range("D44").Value = 102300
$ws.Range("E44").Value = 95200
$ws.Range("D45").Value = 6300
$ws.Range("E45").Value = 6000
$ws.Range("D46").Value = 267000
$ws.Range("E46").Value = 271600
$ws.Range("D47").Value = 900
$ws.Range("E47").Value = 800
$ws.Range("D48").Value = 42000
$ws.Range("E48").Value = 42100
$ws.Range("D49").Value = 100
$ws.Range("E49").Value = 100
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 9100
$ws.Range("E52").Value = 7900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 319000
$ws.Range("E54").Value = 322500
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("D58").Value = 32600
$ws.Range("E58").Value = 36500
$ws.Range("D59").Value = 51200
$ws.Range("E59").Value = 54700
$ws.Range("D60").Value = 83800
$ws.Range("E60").Value = 91200
$ws.Range("D61").Value = 100
$ws.Range("E61").Value = 100
$ws.Range("D62").Value = 12000
$ws.Range("E62").Value = 11300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 167300
$ws.Range("E66").Value = 171400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 55800
$ws.Range("E72").Value = 55000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 151700
$ws.Range("E76").Value = 151200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("E80").Value = 43281
$ws.Range("D81").Value = 700
$ws.Range("E81").Value = 1800
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 100
$ws.Range("E89").Value = 4600
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -700
$ws.Range("E94").Value = -1300
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -4300
$ws.Range("E100").Value = -100
$ws.Range("D101").Value = 1100
$ws.Range("E101").Value = -4800
$ws.Range("D102").Value = -3800
$ws.Range("E102").Value = -1600
